$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q3" right after "总计" (i.e. before the
#    current second sheet, "2022-Q2"). This shifts every following sheet
#    one position to the right, exactly like the target diff expects.
# ---------------------------------------------------------------------
$totalSheet   = $wb.Worksheets.Item(1)
$q2SheetBefore = $wb.Worksheets.Item(2)

$q3Sheet = $wb.Worksheets.Add($q2SheetBefore)
$q3Sheet.Name = "2022-Q3"

# NOTE: after Worksheets.Add(Before:=...) the variable that used to
# reference the "Before" sheet can no longer be trusted - fetch a fresh
# handle to the "2022-Q2" sheet by name before touching it again.
$q2Sheet = $wb.Worksheets.Item("2022-Q2")

# ---------------------------------------------------------------------
# 2. Populate the new "2022-Q3" sheet with the fund-holding detail rows.
#    Re-use the formatting that already exists on the "2022-Q2" sheet
#    (same column layout / header style) so that styles line up with
#    the rest of the workbook.
# ---------------------------------------------------------------------
$q2Sheet.Range("B1:H1").Copy($q3Sheet.Range("B1:H1"))
$q2Sheet.Range("A2").Copy($q3Sheet.Range("A2"))
$q2Sheet.Range("A2").Copy($q3Sheet.Range("A3"))
$q2Sheet.Range("A2").Copy($q3Sheet.Range("A4"))

$q3Sheet.Range("A2").Value = 0
$q3Sheet.Range("A3").Value = 1
$q3Sheet.Range("A4").Value = 2

$q3Sheet.Range("B2").Value = "'002291"
$q3Sheet.Range("C2").Value = "诺安安鑫灵活配置混合"
$q3Sheet.Range("D2").Value = "'2.66"
$q3Sheet.Range("E2").Value = "'77.38"
$q3Sheet.Range("F2").Value = "'4.33"
$q3Sheet.Range("G2").Value = "'0.1152"
$q3Sheet.Range("H2").Value = 3

$q3Sheet.Range("B3").Value = "'002415"
$q3Sheet.Range("C3").Value = "融通通盈灵活配置混合"
$q3Sheet.Range("D3").Value = "'0.25"
$q3Sheet.Range("E3").Value = "'58.75"
$q3Sheet.Range("F3").Value = "'6.07"
$q3Sheet.Range("G3").Value = "'0.0152"
$q3Sheet.Range("H3").Value = 1

$q3Sheet.Range("B4").Value = "'620004"
$q3Sheet.Range("C4").Value = "金元顺安价值增长混合"
$q3Sheet.Range("D4").Value = "'0.34"
$q3Sheet.Range("E4").Value = "'74.99"
$q3Sheet.Range("F4").Value = "'2.00"
$q3Sheet.Range("G4").Value = "'0.0068"
$q3Sheet.Range("H4").Value = 3

# Remove the "number stored as text" quote-prefix styling so the cells
# come out with the default style, matching the rest of the workbook.
$q3Sheet.Range("B2:G4").Style = "Normal"

# ---------------------------------------------------------------------
# 3. Insert a new summary row in the "总计" sheet for the 2022-Q3 quarter,
#    right above the existing 2022-Q2 row, pushing all the others down.
# ---------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("A3").Copy($totalSheet.Range("A2"))

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2:D2").ClearFormats()
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 3
$totalSheet.Range("D2").Value = 0.14

# The "row index" column (A) is a simple 0-based counter; renumber it
# for every row now that a new row has been inserted at the top.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
